$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"   = -12.1272
    "C7"   = -12.9755
    "E7"   = 15.39170000000001
    "E15"  = 16.2416
    "C16"  = -14.018
    "E21"  = 17.11819999999999
    "E22"  = 17.1001
    "E23"  = 16.13089999999998
    "C28"  = -12.5239
    "C29"  = -11.13940000000001
    "C32"  = -13.34150000000001
    "E34"  = 17.37710000000001
    "C40"  = -12.8739
    "E43"  = 17.56720000000001
    "E45"  = 16.43799999999999
    "E50"  = 16.32049999999999
    "E51"  = 17.2443
    "C52"  = -11.1536
    "C57"  = -13.73309999999999
    "C66"  = -11.3441
    "E66"  = 17.23860000000001
    "E67"  = 17.15600000000002
    "E79"  = 18.19290000000002
    "E84"  = 16.43369999999999
    "E92"  = 18.35810000000002
    "E97"  = 16.45809999999999
    "C100" = -12.5334
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
